# Refresh model output (MACHINE LEARNING MODEL V1.0): new game slate + recomputed feature/pred columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (49 updated cells)
$rowUpdates = [ordered]@{
    4 = 230.5   # D2
    5 = -9   # E2
    6 = "Denver"   # F2
    7 = "Brooklyn"   # G2
    8 = 0.5454545454545454   # H2
    9 = 0.5454545454545454   # I2
    10 = 116.8805970149254   # J2
    11 = 113.7462686567164   # K2
    12 = 97.88358208955219   # L2
    13 = 97.95074626865674   # M2
    14 = 118.9134328358209   # N2
    15 = 115.8492537313433   # O2
    16 = 114.9223880597015   # P2
    17 = 114.8880597014925   # Q2
    18 = 77.12388059701497   # R2
    19 = 74.20447761194028   # S2
    20 = 0.363134328358209   # T2
    21 = 0.3902388059701493   # U2
    22 = 0.6079552238805968   # V2
    23 = 0.6038656716417912   # W2
    24 = 0.2607014925373135   # X2
    25 = 0.2556865671641791   # Y2
    26 = 12.70149253731343   # Z2
    27 = 12.22388059701493   # AA2
    28 = 11.82686567164179   # AB2
    29 = 11.65074626865672   # AC2
    30 = 0.2015149253731344   # AD2
    31 = 0.2110597014925373   # AE2
    32 = 1.020791240305025   # AF2
    33 = 0.9934171935084404   # AG2
    34 = 0.9525390541863533   # AH2
    35 = 1.04032716616367   # AI2
    36 = 11.12416005624721   # AJ2
    37 = 11.25331564218959   # AK2
    38 = 0.6865671641791045   # AL2
    39 = 0.5671641791044776   # AM2
    40 = 51.5   # AN2
    41 = 50.5   # AO2
    42 = 75.09999999999999   # AP2
    43 = 76.5   # AQ2
    44 = 0.75   # AR2
    45 = 0.6   # AS2
    46 = 0.4884292375388056   # AT2
    47 = 0.5079099556161815   # AU2
    48 = 2.16   # AV2
    49 = 0.82   # AW2
    50 = 0.02916223103971095   # AX2
    52 = 0.03629353233830845   # AZ2
    54 = 0.5279092519844266   # BB2
}
foreach ($col in $rowUpdates.Keys) {
    $ws.Cells.Item(2, $col).Value = $rowUpdates[$col]
}

# Row 3 (49 updated cells)
$rowUpdates = [ordered]@{
    4 = 220   # D3
    5 = 8   # E3
    6 = "Charlotte"   # F3
    7 = "Cleveland"   # G3
    8 = 0.4477611940298508   # H3
    9 = 0.5454545454545454   # I3
    10 = 111.8235294117647   # J3
    11 = 111.9855072463768   # K3
    12 = 100.4867647058824   # L3
    13 = 94.70434782608693   # M3
    14 = 110.375   # N3
    15 = 116.8666666666667   # O3
    16 = 116.0191176470588   # P3
    17 = 111.2985507246377   # Q3
    18 = 75.62352941176475   # R3
    19 = 77.39855072463766   # S3
    20 = 0.3563382352941177   # T3
    21 = 0.374623188405797   # U3
    22 = 0.5497647058823528   # V3
    23 = 0.5926521739130435   # W3
    24 = 0.2632647058823531   # X3
    25 = 0.2695507246376813   # Y3
    26 = 11.5764705882353   # Z3
    27 = 11.89855072463768   # AA3
    28 = 12.41176470588235   # AB3
    29 = 13.45942028985507   # AC3
    30 = 0.2055294117647059   # AD3
    31 = 0.2096304347826087   # AE3
    32 = 0.9766247110197791   # AF3
    33 = 0.9780393645971774   # AG3
    34 = 0.9270559354725583   # AH3
    35 = 1.003105991976187   # AI3
    36 = 11.55032616864459   # AJ3
    37 = 10.60009947126309   # AK3
    38 = 0.3235294117647059   # AL3
    39 = 0.6086956521739131   # AM3
    40 = 34.5   # AN3
    41 = 46.5   # AO3
    42 = 75.59999999999999   # AP3
    43 = 77.2   # AQ3
    44 = 0.4666666666666667   # AR3
    45 = 0.3529411764705883   # AS3
    46 = 0.5026769579206853   # AT3
    47 = 0.492151092120753   # AU3
    48 = 0.53   # AV3
    49 = 2.04   # AW3
    51 = 0   # AY3
    53 = 0   # BA3
    54 = 0.5706285307576093   # BB3
}
foreach ($col in $rowUpdates.Keys) {
    $ws.Cells.Item(3, $col).Value = $rowUpdates[$col]
}

# Row 4 (51 updated cells)
$rowUpdates = [ordered]@{
    4 = 227.5   # D4
    5 = -7.5   # E4
    6 = "Philadelphia"   # F4
    7 = "Washington"   # G4
    8 = 0.5692307692307692   # H4
    9 = 0.484375   # I4
    10 = 115.0909090909091   # J4
    11 = 113.2238805970149   # K4
    12 = 96.45454545454545   # L4
    13 = 97.78656716417913   # M4
    14 = 118.4818181818182   # N4
    15 = 115.4029850746269   # O4
    16 = 114.359090909091   # P4
    17 = 115.6686567164179   # Q4
    18 = 77.52727272727272   # R4
    19 = 76.34477611940301   # S4
    20 = 0.3911969696969697   # T4
    21 = 0.3629701492537314   # U4
    22 = 0.6071060606060608   # V4
    23 = 0.589373134328358   # W4
    24 = 0.3085909090909091   # X4
    25 = 0.2648805970149253   # Y4
    26 = 11.81818181818182   # Z4
    27 = 12.40746268656716   # AA4
    28 = 12.56212121212121   # AB4
    29 = 10.64029850746269   # AC4
    30 = 0.2396818181818181   # AD4
    31 = 0.2038059701492537   # AE4
    32 = 1.005160778086542   # AF4
    33 = 0.9888548523756762   # AG4
    34 = 1.11216429699842   # AH4
    35 = 1.018630811143334   # AI4
    36 = 12.31250568094288   # AJ4
    37 = 10.8940910938849   # AK4
    38 = 0.6666666666666666   # AL4
    39 = 0.4626865671641791   # AM4
    40 = 50.5   # AN4
    41 = 35.5   # AO4
    42 = 76.2   # AP4
    43 = 76.09999999999999   # AQ4
    44 = 0.7   # AR4
    45 = 0.5555555555555556   # AS4
    46 = 0.4820341586818817   # AT4
    47 = 0.489537703351612   # AU4
    48 = 2.57   # AV4
    49 = -3.23   # AW4
    50 = 0   # AX4
    51 = 0   # AY4
    52 = 0   # AZ4
    53 = 0   # BA4
    54 = 0.5007639899963852   # BB4
}
foreach ($col in $rowUpdates.Keys) {
    $ws.Cells.Item(4, $col).Value = $rowUpdates[$col]
}

# Row 5 (50 updated cells)
$rowUpdates = [ordered]@{
    4 = 230.5   # D5
    5 = 1   # E5
    6 = "NewOrleans"   # F5
    7 = "Portland"   # G5
    8 = 0.4545454545454545   # H5
    9 = 0.5074626865671642   # I5
    10 = 114.0757575757576   # J5
    11 = 114.6268656716418   # K5
    12 = 98.60303030303029   # L5
    13 = 97.2014925373134   # M5
    14 = 114.6636363636364   # N5
    15 = 117.5970149253731   # O5
    16 = 114.0530303030303   # P5
    17 = 118.7835820895522   # Q5
    18 = 77.1878787878788   # R5
    19 = 75.67014925373135   # S5
    20 = 0.3455303030303031   # T5
    21 = 0.4193134328358208   # U5
    22 = 0.5800151515151516   # V5
    23 = 0.6003432835820898   # W5
    24 = 0.2836363636363636   # X5
    25 = 0.3039701492537314   # Y5
    26 = 12.56060606060606   # Z5
    27 = 12.57910447761194   # AA5
    28 = 12.90606060606061   # AB5
    29 = 11.51194029850746   # AC5
    30 = 0.2182651515151515   # AD5
    31 = 0.2306940298507463   # AE5
    32 = 0.9962948259891492   # AF5
    33 = 1.001107997132243   # AG5
    34 = 0.9350511356089787   # AH5
    35 = 0.9363715277777777   # AI5
    36 = 10.44075726683016   # AJ5
    37 = 12.50651868342949   # AK5
    38 = 0.4848484848484849   # AL5
    39 = 0.4626865671641791   # AM5
    40 = 45.5   # AN5
    41 = 39.5   # AO5
    42 = 77.40000000000001   # AP5
    44 = 0.4705882352941176   # AR5
    45 = 0.5   # AS5
    46 = 0.5188261324245436   # AT5
    47 = 0.4870639960032235   # AU5
    48 = -0.84   # AV5
    49 = -8.550000000000001   # AW5
    50 = 0.2500043071515946   # AX5
    51 = 0.03251883489450622   # AY5
    52 = 0.2133217592592592   # AZ5
    53 = 0.04824523217247097   # BA5
    54 = 0.5179657086019411   # BB5
}
foreach ($col in $rowUpdates.Keys) {
    $ws.Cells.Item(5, $col).Value = $rowUpdates[$col]
}

# Row 6 (51 updated cells)
$rowUpdates = [ordered]@{
    4 = 237.5   # D6
    5 = 3.5   # E6
    6 = "SanAntonio"   # F6
    7 = "OklahomaCity"   # G6
    8 = 0.4090909090909091   # H6
    9 = 0.609375   # I6
    10 = 112.3030303030303   # J6
    11 = 118.1818181818182   # K6
    12 = 100.6090909090909   # L6
    13 = 100.6257575757576   # M6
    14 = 111.1590909090909   # N6
    15 = 116.2181818181818   # O6
    16 = 120.9909090909091   # P6
    17 = 115.0924242424242   # Q6
    18 = 75.42272727272727   # R6
    19 = 72.85151515151516   # S6
    20 = 0.3376060606060606   # T6
    21 = 0.3678636363636363   # U6
    22 = 0.5577575757575757   # V6
    23 = 0.574681818181818   # W6
    24 = 0.2331969696969697   # X6
    25 = 0.2529848484848485   # Y6
    26 = 12.95909090909091   # Z6
    27 = 11.03484848484849   # AA6
    28 = 11.68787878787879   # AB6
    29 = 14.08787878787879   # AC6
    30 = 0.1899166666666667   # AD6
    31 = 0.2169772727272727   # AE6
    32 = 0.9808124917295222   # AF6
    33 = 1.032155617308456   # AG6
    34 = 1.03291958985429   # AH6
    35 = 1.035128205128205   # AI6
    36 = 10.54136110033076   # AJ6
    37 = 12.64900174634597   # AK6
    38 = 0.2575757575757576   # AL6
    39 = 0.4696969696969697   # AM6
    40 = 22.5   # AN6
    41 = 23.5   # AO6
    42 = 73.2   # AP6
    43 = 74.40000000000001   # AQ6
    44 = 0.4545454545454545   # AR6
    45 = 0.3157894736842105   # AS6
    46 = 0.4919335263662818   # AT6
    47 = 0.5081090765260694   # AU6
    48 = 0.02   # AV6
    49 = -7.9   # AW6
    50 = 0   # AX6
    51 = 0.04072586953120923   # AY6
    52 = 0   # AZ6
    53 = 0.04276830808080808   # BA6
    54 = 0.4828193237373674   # BB6
}
foreach ($col in $rowUpdates.Keys) {
    $ws.Cells.Item(6, $col).Value = $rowUpdates[$col]
}

# Row 7 (44 updated cells)
$rowUpdates = [ordered]@{
    4 = 225.5   # D7
    5 = -2.5   # E7
    7 = "NewYork"   # G7
    8 = 0.5   # H7
    9 = 0.5606060606060606   # I7
    10 = 116.7164179104478   # J7
    11 = 115.3382352941177   # K7
    12 = 101.3238805970149   # L7
    13 = 96.04117647058824   # M7
    14 = 114.1686567164179   # N7
    15 = 118.2220588235294   # O7
    16 = 114.7611940298508   # P7
    17 = 115.2808823529411   # Q7
    18 = 76.41492537313432   # R7
    19 = 76.61176470588234   # S7
    20 = 0.3472835820895522   # T7
    21 = 0.3987205882352943   # U7
    22 = 0.5793582089552238   # V7
    23 = 0.5737058823529413   # W7
    24 = 0.2943432835820896   # X7
    25 = 0.2870882352941176   # Y7
    26 = 11.85223880597015   # Z7
    27 = 10.66470588235294   # AA7
    28 = 10.48059701492537   # AB7
    29 = 10.72352941176471   # AC7
    30 = 0.2044850746268657   # AD7
    31 = 0.2163529411764706   # AE7
    32 = 1.019357361663299   # AF7
    33 = 1.007320832263036   # AG7
    34 = 0.9910059676044332   # AH7
    35 = 1.020187853287432   # AI7
    36 = 10.52616635881319   # AJ7
    37 = 11.22839139665168   # AK7
    38 = 0.4925373134328358   # AL7
    39 = 0.5735294117647058   # AM7
    41 = 38.5   # AO7
    43 = 75.8   # AQ7
    45 = 0.6190476190476191   # AS7
    46 = 0.5107478346136662   # AT7
    47 = 0.5084661876733064   # AU7
    49 = -0.01   # AW7
    50 = 0.2579974119485635   # AX7
    52 = 0.1334501105583195   # AZ7
    54 = 0.4371395398434437   # BB7
}
foreach ($col in $rowUpdates.Keys) {
    $ws.Cells.Item(7, $col).Value = $rowUpdates[$col]
}

